# "feat: added health and gas mask pickups, high jump next"
#
# - Adds a new credits row (row 7) for a new asset:
#     pixel_icons_by_oceansdream.png / various-inventory-24-pixel-icon-set / CC-BY 3.0, CC-BY-SA 3.0
# - Converts the "Source" (column B) text in several rows into real
#   hyperlinks (styled with the existing built-in "Hyperlink" cell style),
#   matching the one that row 6 already had.
# - Updates the active selection / window scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 was an empty gap row in the sheet; fill it in directly -- the rows
# below (8, 9, 10) keep their original row numbers, nothing shifts.
$ws.Range("A7").Value = "pixel_icons_by_oceansdream.png"
$ws.Range("B7").Value = "https://opengameart.org/content/various-inventory-24-pixel-icon-set"
$ws.Range("C7").Value = "CC-BY 3.0, CC-BY-SA 3.0"

# Add hyperlinks to the Source cells (in the same order the workbook's
# relationship ids were minted: B4, B2, B3, B5, B7, B9, B10), then apply the
# built-in "Hyperlink" style so they render/format like the existing B6 link.
$hyperlinkCells = @(
    @{ Cell = "B4";  Url = "https://www.fontspace.com/a-area-kilometer-50-font-f53888" },
    @{ Cell = "B2";  Url = "https://freesound.org/people/Whiprealgood/sounds/87535/" },
    @{ Cell = "B3";  Url = "https://freesound.org/people/suntemple/sounds/253172/" },
    @{ Cell = "B5";  Url = "https://opengameart.org/content/simple-explosion-bleeds-game-art" },
    @{ Cell = "B7";  Url = "https://opengameart.org/content/various-inventory-24-pixel-icon-set" },
    @{ Cell = "B9";  Url = "https://elthen.itch.io/2d-pixel-art-vegetable-monsters-sprite-pack" },
    @{ Cell = "B10"; Url = "https://free-game-assets.itch.io/night-city-street-2d-background-tiles" }
)

foreach ($link in $hyperlinkCells) {
    $range = $ws.Range($link.Cell)
    $ws.Hyperlinks.Add($range, $link.Url) | Out-Null
    $range.Style = "Hyperlink"
}

# Move the active selection to where the author left off editing, and
# restore the workbook window's on-screen position.
$ws.Range("C14").Select()
$excel.ActiveWindow.Left = 39300
$excel.ActiveWindow.Top = 1880
